# Updated cryptos list — price (D) and 1h volume % (E) refresh, plus a
# rank swap between FraxShare and WEMIXTOKEN (rows 42/43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces text entry so Excel does not reinterpret numeric-
# looking strings (e.g. "7.390", "28.155.60") as numbers and silently
# drop trailing zeros / treat multi-dot values as something else.

$ws.Range("D2").Value = "'28.155.60"
$ws.Range("E2").Value = "'  +1.51%  "
$ws.Range("D3").Value = "'1.803.41"
$ws.Range("E3").Value = "'  +2.57%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'324.64"
$ws.Range("E5").Value = "'  -0.29%  "
$ws.Range("D7").Value = "'0.4284"
$ws.Range("E7").Value = "'  -3.36%  "
$ws.Range("D8").Value = "'0.3634"
$ws.Range("E8").Value = "'  -2.73%  "
$ws.Range("D9").Value = "'44.84"
$ws.Range("E9").Value = "'  -1.32%  "
$ws.Range("D10").Value = "'0.07572"
$ws.Range("E10").Value = "'  +0.48%  "
$ws.Range("D11").Value = "'1.125"
$ws.Range("E11").Value = "'  -0.10%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "'  +0.03%  "
$ws.Range("D13").Value = "'21.66"
$ws.Range("E13").Value = "'  -0.43%  "
$ws.Range("D14").Value = "'6.223"
$ws.Range("E14").Value = "'  +0.22%  "
$ws.Range("D15").Value = "'7.390"
$ws.Range("E15").Value = "'  -0.13%  "
$ws.Range("D16").Value = "'1.823.56"
$ws.Range("E16").Value = "'  +3.66%  "
$ws.Range("D17").Value = "'93.26"
$ws.Range("E17").Value = "'  +5.56%  "
$ws.Range("D18").Value = "'0.00001071"
$ws.Range("E18").Value = "'  -0.16%  "
$ws.Range("D19").Value = "'0.06370"
$ws.Range("E19").Value = "'  +2.53%  "
$ws.Range("E20").Value = "'  +0.12%  "
$ws.Range("D21").Value = "'17.27"
$ws.Range("E21").Value = "'  -0.43%  "
$ws.Range("D22").Value = "'6.096"
$ws.Range("E22").Value = "'  -1.43%  "
$ws.Range("D23").Value = "'28.149.48"
$ws.Range("E23").Value = "'  +1.37%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "'  -1.76%  "
$ws.Range("D25").Value = "'2.150"
$ws.Range("E25").Value = "'  -7.11%  "
$ws.Range("D26").Value = "'160.32"
$ws.Range("E26").Value = "'  +4.63%  "
$ws.Range("D27").Value = "'20.48"
$ws.Range("E27").Value = "'  -1.00%  "
$ws.Range("D28").Value = "'2.026.80"
$ws.Range("E28").Value = "'  +3.53%  "
$ws.Range("D29").Value = "'2.236"
$ws.Range("E29").Value = "'  -5.67%  "
$ws.Range("D30").Value = "'129.25"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("D31").Value = "'1.182"
$ws.Range("E31").Value = "'  -3.12%  "
$ws.Range("D32").Value = "'5.896"
$ws.Range("E32").Value = "'  +2.40%  "
$ws.Range("E33").Value = "'  -3.42%  "
$ws.Range("D34").Value = "'3.539"
$ws.Range("E34").Value = "'  -2.86%  "
$ws.Range("D35").Value = "'12.87"
$ws.Range("E35").Value = "'  +1.24%  "
$ws.Range("D36").Value = "'0.02370"
$ws.Range("E36").Value = "'  +1.48%  "
$ws.Range("D37").Value = "'5.137"
$ws.Range("E37").Value = "'  +1.06%  "
$ws.Range("D38").Value = "'0.6508"
$ws.Range("E38").Value = "'  +0.19%  "
$ws.Range("D39").Value = "'0.2129"
$ws.Range("E39").Value = "'  -2.13%  "
$ws.Range("D40").Value = "'0.06133"
$ws.Range("E40").Value = "'  -0.18%  "
$ws.Range("D41").Value = "'1.198"
$ws.Range("E41").Value = "'  -0.33%  "
$ws.Range("E44").Value = "'  +0.19%  "
$ws.Range("D45").Value = "'13.61"
$ws.Range("E45").Value = "'  -1.37%  "
$ws.Range("D46").Value = "'0.6019"
$ws.Range("E46").Value = "'  +0.17%  "
$ws.Range("D47").Value = "'3.725"
$ws.Range("E47").Value = "'  -0.78%  "
$ws.Range("D48").Value = "'125.37"
$ws.Range("E48").Value = "'  -0.85%  "
$ws.Range("D49").Value = "'1.993"
$ws.Range("E49").Value = "'  +0.19%  "
$ws.Range("D50").Value = "'1.165"
$ws.Range("E50").Value = "'  +2.47%  "
$ws.Range("D51").Value = "'0.06975"
$ws.Range("E51").Value = "'  +1.00%  "

# Rows 42/43: FraxShare and WEMIXTOKEN swap rank order, with refreshed price/volume.
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.427"
$ws.Range("E42").Value = "'  +0.58%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.994"
$ws.Range("E43").Value = "'  -0.01%  "
